$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 74, shifting rows 74:194 down to 75:195
$ws.Rows.Item(74).Insert()

# Fill in the new row 74 with the new record's data.
$ws.Range("A74").Value = 5
$ws.Range("B74").Value = "Macroferia Regional de Talca"
$ws.Range("C74").Value = "Maule"
$ws.Range("D74").Value = 44495
$ws.Range("E74").Value = 7
$ws.Range("F74").Value = 100114014
$ws.Range("G74").Value = "Betarraga"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 650
$ws.Range("L74").Value = 650
$ws.Range("M74").Value = 650
$ws.Range("N74").Value = "$/paquete 5 unidades"
$ws.Range("O74").Value = "Región del Maule"
$ws.Range("P74").Value = 130
$ws.Range("Q74").Value = 5
$ws.Range("R74").Value = "Hortaliza"
